$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the new, longer question text
$ws.Columns.Item(3).ColumnWidth = 50.6

# Copy formatting (fills/borders/fonts) from the last data row down to the new row
$ws.Range("A17:G17").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)

# New LeetCode entry: "Reverse Substrings Between Each Pair of Parentheses"
$ws.Range("A18").Value = 1190
$ws.Range("B18").Value = "Medium"
$ws.Range("C18").Value = "Reverse Substrings Between Each Pair of Parentheses"
$ws.Range("D18").Value = "http://rb.gy/0h2e8g"
$ws.Range("E18").Value = "Stack"
$ws.Range("F18").Value = "O(n2)"
$ws.Range("G18").Value = "Use the stack cleverly"

# Turn the url cell into a real hyperlink (matches the style used by the other rows)
$ws.Hyperlinks.Add($ws.Range("D18"), "http://rb.gy/0h2e8g")
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# Row height tweaks caused by the re-wrapped text after widening column C
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(18).RowHeight = 16.5

# Leave the selection where the author ended up after entering the new row
[void]$ws.Range("A20").Select()
